$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct status name: "bleu" -> "noir" and update the associated label text
$ws.Range("B2").Value = "noir"
$ws.Range("C2").Value = "pas de résultat postés ni publiés"
